$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.7115070223808289
$ws.Range("B1").Value = 1.498084306716919
$ws.Range("D1").Value = 1.980963945388794
$ws.Range("E1").Value = 1.243313789367676
